$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 15.13816947372388
$ws.Range("C2").Value = 10.57619841282403
$ws.Range("E2").Value = 12.48759964142782
$ws.Range("F2").Value = 16.86991607391233
$ws.Range("G2").Value = 3.640653167366422
$ws.Range("L2").Value = 9.867383632438891
$ws.Range("M2").Value = 14.84076601712459
$ws.Range("N2").Value = 18.29681617104089
$ws.Range("O2").Value = 23.09019337217143

$ws.Range("B3").Value = 14.73060681304302
$ws.Range("C3").Value = 10.43455430430983
$ws.Range("E3").Value = 12.52719117488554
$ws.Range("F3").Value = 15.89584955866808
$ws.Range("G3").Value = 3.64270253703804
$ws.Range("L3").Value = 9.874631915804647
$ws.Range("M3").Value = 14.76268644959621
$ws.Range("N3").Value = 18.35256908831414
$ws.Range("O3").Value = 23.16283643608086

$ws.Range("B4").Value = 14.47698313372717
$ws.Range("C4").Value = 10.34587084369377
$ws.Range("E4").Value = 12.55291072002945
$ws.Range("F4").Value = 15.26997757108491
$ws.Range("G4").Value = 3.644027937752507
$ws.Range("L4").Value = 9.880418706905916
$ws.Range("M4").Value = 14.71681313574361
$ws.Range("N4").Value = 18.38872071293296
$ws.Range("O4").Value = 23.21329246586749

$ws.Range("B5").Value = 14.37293757822897
$ws.Range("C5").Value = 10.30932489346967
$ws.Range("E5").Value = 12.56374700129244
$ws.Range("F5").Value = 15.00819731993403
$ws.Range("G5").Value = 3.644584970366805
$ws.Range("L5").Value = 9.883113408712298
$ws.Range("M5").Value = 14.69865361350314
$ws.Range("N5").Value = 18.40393645153275
$ws.Range("O5").Value = 23.23532066324243

$ws.Range("B6").Value = 14.35562388472703
$ws.Range("C6").Value = 10.3032326171877
$ws.Range("E6").Value = 12.56556784258816
$ws.Range("F6").Value = 14.96433081551593
$ws.Range("G6").Value = 3.644678488769679
$ws.Range("L6").Value = 9.883581201534748
$ws.Range("M6").Value = 14.69567090722551
$ws.Range("N6").Value = 18.40649225718603
$ws.Range("O6").Value = 23.23906688431452

$ws.Range("B7").Value = 14.47558252551432
$ws.Range("C7").Value = 10.34537958673934
$ws.Range("E7").Value = 12.55305542206558
$ws.Range("F7").Value = 15.26647399323137
$ws.Range("G7").Value = 3.644035381505101
$ws.Range("L7").Value = 9.880453685343282
$ws.Range("M7").Value = 14.71656604872224
$ws.Range("N7").Value = 18.38892395795228
$ws.Range("O7").Value = 23.21358361273264

$ws.Range("B8").Value = 14.99843498003424
$ws.Range("C8").Value = 10.52772943558097
$ws.Range("E8").Value = 12.50095857379506
$ws.Range("F8").Value = 16.53996406344768
$ws.Range("G8").Value = 3.641345899879338
$ws.Range("L8").Value = 9.869605812022581
$ws.Range("M8").Value = 14.81342374910923
$ws.Range("N8").Value = 18.31564211620533
$ws.Range("O8").Value = 23.11402311210457

$ws.Range("B9").Value = 15.99047237285147
$ws.Range("C9").Value = 10.87069495838726
$ws.Range("E9").Value = 12.40995152741264
$ws.Range("F9").Value = 19.00274580682531
$ws.Range("G9").Value = 3.636601643489184
$ws.Range("L9").Value = 9.858912081805055
$ws.Range("M9").Value = 15.01911354457428
$ws.Range("N9").Value = 18.18711481060626
$ws.Range("O9").Value = 22.96543094594642

$ws.Range("B10").Value = 16.69106755503867
$ws.Range("C10").Value = 11.11243416168783
$ws.Range("E10").Value = 12.34984067645498
$ws.Range("F10").Value = 20.67494806633232
$ws.Range("G10").Value = 3.633435613230536
$ws.Range("L10").Value = 9.8574695212699
$ws.Range("M10").Value = 15.17891990619752
$ws.Range("N10").Value = 18.10187054606856
$ws.Range("O10").Value = 22.88496382297739

$ws.Range("B11").Value = 17.00217012951837
$ws.Range("C11").Value = 11.21991933819238
$ws.Range("E11").Value = 12.3239506763451
$ws.Range("F11").Value = 21.3917225636224
$ws.Range("G11").Value = 3.632063971237305
$ws.Range("L11").Value = 9.858197148648685
$ws.Range("M11").Value = 15.25329998174002
$ws.Range("N11").Value = 18.06507087012443
$ws.Range("O11").Value = 22.85464159940378

$ws.Range("B12").Value = 17.11877395437608
$ws.Range("C12").Value = 11.26024292062588
$ws.Range("E12").Value = 12.31435524667529
$ws.Range("F12").Value = 21.65686569030329
$ws.Range("G12").Value = 3.631554375862884
$ws.Range("L12").Value = 9.858670780335855
$ws.Range("M12").Value = 15.28168980073434
$ws.Range("N12").Value = 18.05141924483702
$ws.Range("O12").Value = 22.84406655167632

$ws.Range("B13").Value = 17.09371643908727
$ws.Range("C13").Value = 11.25157569415379
$ws.Range("E13").Value = 12.31641252878674
$ws.Range("F13").Value = 21.60004134736742
$ws.Range("G13").Value = 3.631663690566675
$ws.Range("L13").Value = 9.858559980146584
$ws.Range("M13").Value = 15.27556589209342
$ws.Range("N13").Value = 18.05434676649131
$ws.Range("O13").Value = 22.84630366586153

$ws.Range("B14").Value = 17.01178792477402
$ws.Range("C14").Value = 11.22324449440197
$ws.Range("E14").Value = 12.32315707849877
$ws.Range("F14").Value = 21.4136618050453
$ws.Range("G14").Value = 3.632021850077782
$ws.Range("L14").Value = 9.858232149938503
$ws.Range("M14").Value = 15.25563125234554
$ws.Range("N14").Value = 18.06394206358631
$ws.Range("O14").Value = 22.85375337877184

$ws.Range("B15").Value = 16.96144443437552
$ws.Range("C15").Value = 11.20584086927829
$ws.Range("E15").Value = 12.32731545078266
$ws.Range("F15").Value = 21.29868154950795
$ws.Range("G15").Value = 3.632242509961854
$ws.Range("L15").Value = 9.858057114180678
$ws.Range("M15").Value = 15.24344929132929
$ws.Range("N15").Value = 18.06985636940108
$ws.Range("O15").Value = 22.85843481240178

$ws.Range("B16").Value = 16.6705735969925
$ws.Range("C16").Value = 11.10535791356974
$ws.Range("E16").Value = 12.35156186704146
$ws.Range("F16").Value = 20.62722412089977
$ws.Range("G16").Value = 3.63352662948554
$ws.Range("L16").Value = 9.857449740531976
$ws.Range("M16").Value = 15.17409130462694
$ws.Range("N16").Value = 18.10431522606614
$ws.Range("O16").Value = 22.88707221418882

$ws.Range("B17").Value = 16.49010468203814
$ws.Range("C17").Value = 11.04306318213597
$ws.Range("E17").Value = 12.36680839255883
$ws.Range("F17").Value = 20.20408069597325
$ws.Range("G17").Value = 3.634331930560651
$ws.Range("L17").Value = 9.857430958413255
$ws.Range("M17").Value = 15.13196061874101
$ws.Range("N17").Value = 18.1259606972748
$ws.Range("O17").Value = 22.90625234020821

$ws.Range("B18").Value = 16.38559420196506
$ws.Range("C18").Value = 11.00700040721842
$ws.Range("E18").Value = 12.37571475405049
$ws.Range("F18").Value = 19.95656407809801
$ws.Range("G18").Value = 3.63480157840313
$ws.Range("L18").Value = 9.85755046937536
$ws.Range("M18").Value = 15.10788773146187
$ws.Range("N18").Value = 18.13859686389948
$ws.Range("O18").Value = 22.91787538468753

$ws.Range("B19").Value = 16.35009043032512
$ws.Range("C19").Value = 10.99475091467367
$ws.Range("E19").Value = 12.37875383815204
$ws.Range("F19").Value = 19.87204792380568
$ws.Range("G19").Value = 3.634961704099728
$ws.Range("L19").Value = 9.857613341465708
$ws.Range("M19").Value = 15.09976502553032
$ws.Range("N19").Value = 18.14290727083264
$ws.Range("O19").Value = 22.92191215373101

$ws.Range("B20").Value = 16.50939017529545
$ws.Range("C20").Value = 11.04971878446295
$ws.Range("E20").Value = 12.36517120335118
$ws.Range("F20").Value = 20.24955283636154
$ws.Range("G20").Value = 3.634245536661027
$ws.Range("L20").Value = 9.857419477081788
$ws.Range("M20").Value = 15.13642911934964
$ws.Range("N20").Value = 18.12363722866749
$ws.Range("O20").Value = 22.90414937276975

$ws.Range("B21").Value = 17.03588582417542
$ws.Range("C21").Value = 11.23157650080587
$ws.Range("E21").Value = 12.32117038575414
$ws.Range("F21").Value = 21.46857628470577
$ws.Range("G21").Value = 3.631916383942183
$ws.Range("L21").Value = 9.858323072952819
$ws.Range("M21").Value = 15.26148061613915
$ws.Range("N21").Value = 18.06111600190134
$ws.Range("O21").Value = 22.85154056382082

$ws.Range("B22").Value = 17.37291914458892
$ws.Range("C22").Value = 11.34821417238373
$ws.Range("E22").Value = 12.29362861728482
$ws.Range("F22").Value = 22.22866616901552
$ws.Range("G22").Value = 3.63045133879664
$ws.Range("L22").Value = 9.86006780035445
$ws.Range("M22").Value = 15.34450415024416
$ws.Range("N22").Value = 18.02190749212456
$ws.Range("O22").Value = 22.82244778541717

$ws.Range("B23").Value = 17.19371832540787
$ws.Range("C23").Value = 11.28617220449334
$ws.Range("E23").Value = 12.30821718425494
$ws.Range("F23").Value = 21.82633154458858
$ws.Range("G23").Value = 3.631228044112775
$ws.Range("L23").Value = 9.859031317641213
$ws.Range("M23").Value = 15.30008056983164
$ws.Range("N23").Value = 18.04268286994781
$ws.Range("O23").Value = 22.83748990722527

$ws.Range("B24").Value = 16.500673551561
$ws.Range("C24").Value = 11.04671055968953
$ws.Range("E24").Value = 12.36591093756227
$ws.Range("F24").Value = 20.22900810905287
$ws.Range("G24").Value = 3.634284574563924
$ws.Range("L24").Value = 9.857424261826743
$ws.Range("M24").Value = 15.13440844592823
$ws.Range("N24").Value = 18.12468707099623
$ws.Range("O24").Value = 22.90509826769756

$ws.Range("B25").Value = 15.72652462265632
$ws.Range("C25").Value = 10.77961861600232
$ws.Range("E25").Value = 12.43338217370105
$ws.Range("F25").Value = 18.34778573295695
$ws.Range("G25").Value = 3.637828724864514
$ws.Range("L25").Value = 9.860675800305781
$ws.Range("M25").Value = 14.96187506209794
$ws.Range("N25").Value = 18.2202670651809
$ws.Range("O25").Value = 23.00060666929713

Write-Host "Updated loading_percent values for rows 2-25"